# Add a hydrogen export process
$wb = $excel.ActiveWorkbook

$wsTech = $wb.Worksheets.Item("tech_data")
$wsProc = $wb.Worksheets.Item("processes")

# --- processes sheet: add new IRE process row (row 15) ---
$wsProc.Range("B15").Value = "IRE"
$wsProc.Range("C15").Value = "EXPH2GMD"
$wsProc.Range("D15").Value = "Hydrogen (medium pressure gas) - export"
$wsProc.Range("E15").Value = "PJ"
$wsProc.Range("F15").Value = "PJa"
$wsProc.Range("G15").Value = "annual"

# --- tech_data sheet: add new header column Q (IRE_PRICE) ---
$wsTech.Range("P3").Copy()
$wsTech.Range("Q3").PasteSpecial(-4122)
$wsTech.Range("Q3").Value = "IRE_PRICE"

# --- tech_data sheet: add new data row 25 for the export process ---
$wsTech.Range("B25").Formula = "=processes!C15"
$wsTech.Range("C25").Formula = "=processes!D15"
$wsTech.Range("D25").Value = "H2GMD"
$wsTech.Range("F25").Value = "exp"
$wsTech.Range("G25").Value = 2023
$wsTech.Range("H25").Value = 2023
$wsTech.Range("Q25").Formula = "=5*0.12"

# --- tech_data sheet: convert the NCAP_FOM (col O) formulas for rows 17:24 ---
# into a single shared-formula group, as Excel would when re-entering an
# identical formula across the block.
$wsTech.Range("O17:O24").Formula = "=N17*0.02"

# Re-entering the formulas above resets formatting on previously-unstyled
# cells (Excel's auto-fill formatting behaviour); restore those cells back
# to the Normal style to match their original (unstyled) appearance.
$wsTech.Range("O19").Style = "Normal"
$wsTech.Range("O21").Style = "Normal"
$wsTech.Range("O23").Style = "Normal"
$wsTech.Range("O24").Style = "Normal"

# --- update the selection markers left behind by the editing session ---
# (select the "processes" sheet first, then "tech_data" last so that
# "tech_data" remains the active/selected tab, as in the target workbook)
$wsProc.Range("D23").Select()
$wsTech.Range("P25").Select()
